$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# --- Cells that change numeric <-> text representation: copy formatting from a
# donor cell that already carries the target style/shared-string, then (for the
# numeric targets) overwrite the value. This preserves the existing style index
# instead of Excel fabricating a brand-new one. ---

# Row 16: C16 goes from blank-text "0" to the number 4
$ws.Range("D16").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 4

# Row 18: C18 goes from the number 2 to blank-text "0"
$ws.Range("C14").Copy($ws.Range("C18"))

# Row 26: D26 (was 1) and E26 (was -100) both become blank-text markers
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))

# Row 27: C27 goes from the number 1 to blank-text "0"
$ws.Range("C14").Copy($ws.Range("C27"))

# Row 28: C28 goes from the number 1 to blank-text "0"
$ws.Range("C14").Copy($ws.Range("C28"))

# Row 29: C29 goes from the number 1 to blank-text "0"
$ws.Range("C14").Copy($ws.Range("C29"))

# --- Row 15 (Rape) ---
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = -67.692307692307

# --- Row 16 (Robbery) ---
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 125
$ws.Range("I16").Value = 114
$ws.Range("J16").Value = 104
$ws.Range("K16").Value = 9.615384615384
$ws.Range("L16").Value = -5.785123966942
$ws.Range("M16").Value = -58.844765342960
$ws.Range("N16").Value = -87.513691128149

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 285
$ws.Range("J17").Value = 321
$ws.Range("K17").Value = -11.214953271028
$ws.Range("L17").Value = -11.490683229813
$ws.Range("M17").Value = -0.349650349650
$ws.Range("N17").Value = -52.814569536423

# --- Row 18 (Burglary) ---
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 103
$ws.Range("K18").Value = -38.834951456310
$ws.Range("L18").Value = -39.423076923076
$ws.Range("M18").Value = -79.611650485436
$ws.Range("N18").Value = -92.953020134228

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -53.846153846153
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 327
$ws.Range("J19").Value = 361
$ws.Range("K19").Value = -9.418282548476
$ws.Range("L19").Value = 18.478260869565
$ws.Range("M19").Value = -36.1328125
$ws.Range("N19").Value = -90.227136879856

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 16
$ws.Range("H20").Value = 45.454545454545
$ws.Range("I20").Value = 171
$ws.Range("J20").Value = 192
$ws.Range("K20").Value = -10.9375
$ws.Range("L20").Value = 19.580419580419
$ws.Range("M20").Value = -22.272727272727
$ws.Range("N20").Value = -87.820512820512

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 28
$ws.Range("E21").Value = -35.714285714285
$ws.Range("F21").Value = 80
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -3.614457831325
$ws.Range("I21").Value = 984
$ws.Range("J21").Value = 1122
$ws.Range("K21").Value = -12.299465240641
$ws.Range("L21").Value = -1.303911735205
$ws.Range("M21").Value = -40
$ws.Range("N21").Value = -86.442546155965

# --- Row 23 (Housing) ---
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("I23").Value = 13
$ws.Range("K23").Value = 225
$ws.Range("L23").Value = -23.529411764705
$ws.Range("M23").Value = 30

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -21.428571428571
$ws.Range("F24").Value = 71
$ws.Range("G24").Value = 124
$ws.Range("H24").Value = -42.741935483871
$ws.Range("I24").Value = 937
$ws.Range("J24").Value = 1080
$ws.Range("K24").Value = -13.240740740740
$ws.Range("L24").Value = 22.643979057591
$ws.Range("M24").Value = 10.365135453474

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -11.111111111111
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 30
$ws.Range("H25").Value = 53.333333333333
$ws.Range("I25").Value = 539
$ws.Range("J25").Value = 421
$ws.Range("K25").Value = 28.028503562945
$ws.Range("L25").Value = 35.427135678392
$ws.Range("M25").Value = -18.209408194233

# --- Row 26 (UCR Rape*) ---
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 100
$ws.Range("L26").Value = -2.631578947368

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 47
$ws.Range("K27").Value = -17.021276595744
$ws.Range("L27").Value = -7.142857142857

# --- Row 28 (Shooting Vic.) ---
$ws.Range("N28").Value = -86.71875

# --- Row 29 (Shooting Inc.) ---
$ws.Range("N29").Value = -89.473684210526
